$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.288.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.02%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.834.97'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.022'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +1.70%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.72'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.017'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4332'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3713'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.38%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07318'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.72%  '

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8749'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.76%  '

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2.095.30'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +17.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.30'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.470'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.669'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07095'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.91%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.94'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.23%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009000'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.015'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.40'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.306.35'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.222'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.08'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.348.22'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +17.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.76'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.897'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.53'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.269'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.922'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +7.44%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.41'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08999'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.199'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7594'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.456'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.833'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.018'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.144'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05254'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01946'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5152'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.778'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1660'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.523'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.441'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.89%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '107.81'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.56%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.47'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.24%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.019'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.52%  '

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.913'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4625'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.99%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.661'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.11%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06269'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.01%  '
